$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "rother"
$ws.Range("D2").Value = 784512
$ws.Range("E2").Value = "admin"

$ws.Range("E2").Select()
